$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cells whose new value is plain text (inc. volume %, links, names, and
#     price strings that are not valid plain numbers e.g. "37.709.55") ---
$ws.Range('D2').Value = '37.709.55'
$ws.Range('E2').Value = '  -1.55%  '
$ws.Range('D3').Value = '2.027.81'
$ws.Range('E3').Value = '  -1.97%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('E6').Value = '  -2.03%  '
$ws.Range('E7').Value = '  -2.89%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -3.63%  '
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('D12').Value = '2.327.11'
$ws.Range('E12').Value = '  -2.01%  '
$ws.Range('E13').Value = '  -3.91%  '
$ws.Range('E14').Value = '  -2.29%  '
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('E16').Value = '  -2.92%  '
$ws.Range('D17').Value = '2.025.54'
$ws.Range('E17').Value = '  -2.40%  '
$ws.Range('D18').Value = '37.652.39'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('E19').Value = '  -1.34%  '
$ws.Range('E20').Value = '  -6.90%  '
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('E22').Value = '  -2.06%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('E24').Value = '  -1.67%  '
$ws.Range('E25').Value = '  -1.12%  '
$ws.Range('E26').Value = '  +0.72%  '
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('E28').Value = '  -3.83%  '
$ws.Range('E29').Value = '  -2.22%  '
$ws.Range('E30').Value = '  -6.13%  '
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('E32').Value = '  +7.48%  '
$ws.Range('E33').Value = '  -4.72%  '
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('E35').Value = '  -4.16%  '
$ws.Range('E36').Value = '  +1.93%  '
$ws.Range('E37').Value = '  -2.02%  '
$ws.Range('E38').Value = '  +1.52%  '
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('E40').Value = '  +2.70%  '
$ws.Range('D41').Value = '1.533.19'
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('E42').Value = '  -1.92%  '
$ws.Range('E43').Value = '  -3.14%  '
$ws.Range('E44').Value = '  -2.32%  '
$ws.Range('E45').Value = '  -2.67%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E47').Value = '  -2.96%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E48').Value = '  -2.40%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').Value = '2.215.79'
$ws.Range('E51').Value = '  -2.02%  '

# --- Update Price cells whose new value LOOKS like a plain number (e.g. "21.03").
#     Excel auto-coerces such a string to a Double on assignment, which would lose
#     the original text formatting (trailing zeros, etc). Force the cell to Text
#     first, assign, then restore the default "Normal" style so no stray number
#     format is left behind on the cell.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.607'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.82'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.375'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.42'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.769'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.41'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.89'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '222.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.37'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.89'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.77'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.26'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.20'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0604'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.30'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.40'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.87'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '95.54'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0907'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.06'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.96'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.10'
$ws.Range('D50').Style = 'Normal'
